# Update "想去人数" (want-to-go count) values that changed between crawls.
# Sheet "展览" (Exhibitions) — rows 3,5,7,10,14,18 in column F
# Sheet "全部类型" (All types)   — rows 4,6,8,11,15,19 in column F

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value  = 12712
$wsExhibition.Range("F5").Value  = 77
$wsExhibition.Range("F7").Value  = 33
$wsExhibition.Range("F10").Value = 12612
$wsExhibition.Range("F14").Value = 5878
$wsExhibition.Range("F18").Value = 115

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value  = 12712
$wsAllTypes.Range("F6").Value  = 77
$wsAllTypes.Range("F8").Value  = 33
$wsAllTypes.Range("F11").Value = 12612
$wsAllTypes.Range("F15").Value = 5878
$wsAllTypes.Range("F19").Value = 115
